$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin/Link/Price/Volume(1h) updates from the Sep 4 2023 GitHub Actions refresh.
# D-column "Price" cells are forced to Text (matching the source feed's inline
# string cells) via a NumberFormat round-trip, then ClearFormats() strips the
# leftover style index so the cell ends up unstyled text, same as the original.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.839.04'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.22%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.629.29'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.51%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.54'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.48%  '

$ws.Range("E6").Value = '  +0.08%  '

$ws.Range("E7").Value = '  -0.10%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2569'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.68%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06328'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.61%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.47'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.22%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07759'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.18%  '

$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.239'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.72%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.629.02'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.66%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.855.31'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.47%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5480'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.91%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.874.04'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.25%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.002'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.423'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.55%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '193.73'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.25%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.881'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.25%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.034'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.38%  '

$ws.Range("E24").Value = '  -0.24%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.919'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.88%  '

$ws.Range("E26").Value = '  +0.52%  '

$ws.Range("E27").Value = '  +4.46%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.775'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.11%  '

$ws.Range("E29").Value = '  -0.72%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.238'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.32%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04862'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.38%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.236'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.37%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.183'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.40%  '

$ws.Range("E34").Value = '  +0.57%  '

$ws.Range("E35").Value = '  +0.38%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8946'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.16%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5519'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.96%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.537'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.61%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.122.97'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.84%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01551'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.01%  '

$ws.Range("E41").Value = '  -0.11%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.564'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.14%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7960'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.13'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.27%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.768.30'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.33%  '

$ws.Range("E46").Value = '  -6.39%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4446'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.02%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.003'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.11%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.60'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.26%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05133'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.37%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.524'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.15%  '
